$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6476.811239874914
$ws.Range("C2").Value = 518.105942736684
$ws.Range("D2").Value = 8779.892976837429
